$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constant")

# Update India's proba_infection_per_contact value (row 19, column B)
$ws.Range("B19").Value = 0.005

# Delete row 25 (duplicate "perc_cdr_smearpos" row with the same data as "perc_cdr"),
# shifting all subsequent rows up by one.
$ws.Rows("25:25").Delete()

# Reposition the view to reflect where the edit was made.
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Range("A26").Select()
